# Update "想去人数" (want-to-go count) values on both the "展览" and
# "全部类型" worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 3-6, column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 801
$wsExpo.Range("F4").Value = 273
$wsExpo.Range("F5").Value = 916
$wsExpo.Range("F6").Value = 2212

# Sheet "全部类型": rows 3, 4, 7, 8, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 801
$wsAll.Range("F4").Value = 273
$wsAll.Range("F7").Value = 916
$wsAll.Range("F8").Value = 2212
